$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices stored as text (OOXML inlineStr); some new
# values parse as plain numbers (single decimal point, e.g. "239.81"), so a
# bare Value assignment would silently convert them to numeric cells. Force
# those through Text format while assigning, then restore the default/Normal
# style so the saved cell format matches the originals (unstyled) - only the
# stored text changes, exactly like the diff.

$ws.Range("D2").Value = '26.252.86'
$ws.Range("E2").Value = '  +2.89%  '
$ws.Range("D3").Value = '1.718.38'
$ws.Range("E3").Value = '  +3.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4716'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2623'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06197'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D10").Value = '1.718.43'
$ws.Range("E10").Value = '  +3.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07076'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.32'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5977'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.427'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = '26.269.70'
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006803'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.00%  '
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '1.937.67'
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.538'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.732'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.289'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.763'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.974'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.674'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07763'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04461'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.44%  '
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9758'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6178'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9257'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '113.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.424'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.920'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9997'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01480'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.489'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3824'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1180'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.285'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.781'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3377'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.216'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.58%  '
